# "update apis for getting all departments"
# Sheet1 held a 5-row id/name lookup table (Architect, Technical Lead,
# SW Delivery Manager, Team Lead, Product Manager). It is replaced with a
# 3-row table of department codes (PS-EC, NE-EM, AS); the now-unused rows
# 5 and 6 are cleared back to blank (only their pre-existing style stays).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "PS-EC"
$ws.Range("B3").Value = "NE-EM"
$ws.Range("B4").Value = "AS"

$ws.Range("A5:B6").ClearContents()

# Restore the view: scrolled so row 3 is at the top, with D6 as the active
# selected cell (matches the author's on-screen state when they saved).
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D6").Select()
